$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row 44 with value 0 (matches the pattern of the existing rows)
$ws.Range("A44").Value = 0

# Update the selection to match the recorded UI state after the edit
$ws.Range("B43").Select()
